$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# Row 1 headers for new columns J-Y
$ws.Range("J1").Value = "I4 Low Pixel Limit"
$ws.Range("K1").Value = "I4 High Pixel Limit"
$ws.Range("L1").Value = "I2 Low Pixel Limit"
$ws.Range("M1").Value = "I2 High Pixel Limit"
$ws.Range("N1").Value = "I1 Low Pixel Limit"
$ws.Range("O1").Value = "I1 High Pixel Limit"
$ws.Range("P1").Value = "I3 Low Pixel Limit"
$ws.Range("Q1").Value = "I3 High Pixel Limit"
$ws.Range("R1").Value = "I4 Contrast"
$ws.Range("S1").Value = "I4 Bias"
$ws.Range("T1").Value = "I2 Contrast"
$ws.Range("U1").Value = "I2 Bias"
$ws.Range("V1").Value = "I1 Contrast"
$ws.Range("W1").Value = "I1 Bias"
$ws.Range("X1").Value = "I3 Contrast"
$ws.Range("Y1").Value = "I3 Bias"

# Data rows 2-20 share the same values across all new columns
for ($row = 2; $row -le 20; $row++) {
    $ws.Range("J$row").Value = 2.05987
    $ws.Range("K$row").Value = 14.1319
    $ws.Range("L$row").Value = -0.08746
    $ws.Range("M$row").Value = 14.2581
    $ws.Range("N$row").Value = -0.386169
    $ws.Range("O$row").Value = 14.5226
    $ws.Range("P$row").Value = 1.08658
    $ws.Range("Q$row").Value = 14.124
    $ws.Range("R$row").Value = 1.49254
    $ws.Range("S$row").Value = 0.429104
    $ws.Range("T$row").Value = 1.15672
    $ws.Range("U$row").Value = 0.455224
    $ws.Range("V$row").Value = 0.932836
    $ws.Range("W$row").Value = 0.514925
    $ws.Range("X$row").Value = 1.30597014925
    $ws.Range("Y$row").Value = 0.4738805970149254
}
